$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 116
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 1222
$ws.Range("E2").Value = 1499
$ws.Range("F2").Value = 106
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 6240
$ws.Range("I2").Value = 19322
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 159
$ws.Range("L2").Value = 4440

$ws.Range("B3").Value = 262

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 37
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 2157
$ws.Range("I4").Value = 2162
$ws.Range("J4").Value = 1

$ws.Range("B5").Value = 1358
$ws.Range("D5").Value = 163
$ws.Range("E5").Value = 165
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 650
$ws.Range("I5").Value = 25408
$ws.Range("J5").Value = 4

$ws.Range("B6").Value = 1799
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 342
$ws.Range("E6").Value = 353
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 2197
$ws.Range("I6").Value = 35599
$ws.Range("J6").Value = 9

$ws.Range("B7").Value = 132
$ws.Range("D7").Value = 48
$ws.Range("E7").Value = 53
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 588
$ws.Range("I7").Value = 7917
$ws.Range("J7").Value = 1

$ws.Range("B8").Value = 807
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 875
$ws.Range("E8").Value = 1012
$ws.Range("F8").Value = 61
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 808
$ws.Range("I8").Value = 12063
$ws.Range("J8").Value = 14
$ws.Range("K8").Value = 64
$ws.Range("L8").Value = 841

$ws.Range("B9").Value = 261

$ws.Range("B10").Value = 1020
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 711
$ws.Range("E10").Value = 854
$ws.Range("F10").Value = 136
$ws.Range("G10").Value = 7
$ws.Range("H10").Value = 4832
$ws.Range("I10").Value = 9563
$ws.Range("J10").Value = 4

$ws.Range("B11").Value = 1772
$ws.Range("D11").Value = 212
$ws.Range("E11").Value = 237
$ws.Range("F11").Value = 24
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 2542
$ws.Range("I11").Value = 10263
$ws.Range("J11").Value = 1

$ws.Range("B12").Value = 1386

$ws.Range("B13").Value = 866
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 302
$ws.Range("E13").Value = 305
$ws.Range("F13").Value = 2
$ws.Range("H13").Value = 556
$ws.Range("I13").Value = 12547
$ws.Range("J13").Value = 4

$ws.Range("B14").Value = 18
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 33
$ws.Range("H14").Value = 303

$ws.Range("B15").Value = 469
$ws.Range("D15").Value = 73
$ws.Range("E15").Value = 80
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1017
$ws.Range("I15").Value = 24369

$ws.Range("B16").Value = 186
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = 6667

$ws.Range("B17").Value = 1797
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = 37
$ws.Range("F17").Value = 1
$ws.Range("H17").Value = 417
$ws.Range("I17").Value = 21319

$ws.Range("B18").Value = 3642
$ws.Range("C18").Value = 20
$ws.Range("D18").Value = 477
$ws.Range("E18").Value = 1430
$ws.Range("F18").Value = 116
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 2818
$ws.Range("I18").Value = 16631
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = 802
$ws.Range("L18").Value = 22002
